$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook window geometry ---
$aw = $excel.ActiveWindow
$aw.Top = 0
$aw.Left = 0
$aw.Width = 12760
$aw.Height = 16000

# --- Manager / team data table (rows 2-12: name, dob, nationality) ---
# Row 2 (Hull City FC) - only the manager name changed
$ws.Range("A2").Value = "Marco Silva"

# Row 3 (Leicester City FC) - new data
$ws.Range("A3").Value = "Criag Shakespeare"
$ws.Range("C3").Value = "1963-10-26"
$ws.Range("D3").Value = "England"

# Row 4 (Southampton FC)
$ws.Range("A4").Value = "Claude Puel"
$ws.Range("C4").Value = "1961-09-02"
$ws.Range("D4").Value = "France"

# Row 5 (Watford FC)
$ws.Range("A5").Value = "Walter Mazzarri"
$ws.Range("C5").Value = "1961-10-01"
$ws.Range("D5").Value = "Italy "

# Row 6 (Middlesbrough FC)
$ws.Range("A6").Value = "Steve Agnew"
$ws.Range("C6").Value = "1965-11-09"
$ws.Range("D6").Value = "England"

# Row 7 (Stoke City FC)
$ws.Range("A7").Value = "Mark Hughes"
$ws.Range("C7").Value = "1963-11-01"
$ws.Range("D7").Value = "Wales"

# Row 8 (Everton FC)
$ws.Range("A8").Value = "Ronald Koeman"
$ws.Range("C8").Value = "1963-03-21"
$ws.Range("D8").Value = "Neatherlands"

# Row 9 (Tottenham Hotspur FC) - name rendered with an explicit black font
$ws.Range("A9").Value = "Mauricio Pochettino"
$ws.Range("A9").Font.Color = 0
$ws.Range("C9").Value = "1972-03-02"
$ws.Range("D9").Value = "Argentina "

# Row 10 (Crystal Palace FC)
$ws.Range("A10").Value = "Sam Allardyce"
$ws.Range("C10").Value = "1954-10-19"
$ws.Range("D10").Value = "England"

# Row 11 (West Bromwich Albion FC)
$ws.Range("A11").Value = "Tony Pulis"
$ws.Range("C11").Value = "1958-01-16"
$ws.Range("D11").Value = "Wales"

# Row 12 (Burnley FC)
$ws.Range("A12").Value = "Sean Dyche"
$ws.Range("C12").Value = "1971-06-28"
$ws.Range("D12").Value = "England"

# --- Column widths: A widened for the longer manager names, D sized for the new nationality column ---
$ws.Columns.Item(1).ColumnWidth = 16.83
$ws.Columns.Item(4).ColumnWidth = 11.17

# --- Selection moved to B13 ---
$ws.Range("B13").Select()
